$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.940.37"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "1.875.50"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.22"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5155"
$ws.Range("E7").Value = "  +1.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3719"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07188"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8980"
$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07566"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.89"
$ws.Range("E13").Value = "  +5.03%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.856.45"
$ws.Range("E14").Value = "  -0.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.249"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008490"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("E18").Value = "  +1.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").Value = "26.969.85"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.027"
$ws.Range("E21").Value = "  +0.64%  "

$ws.Range("D22").Value = "2.133.96"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.39"
$ws.Range("E23").Value = "  +1.41%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.94"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.783"
$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.02"
$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.110"
$ws.Range("E28").Value = "  +3.35%  "

$ws.Range("E29").Value = "  +1.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.892"
$ws.Range("E30").Value = "  +5.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.742"
$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09179"
$ws.Range("E32").Value = "  -0.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05034"
$ws.Range("E33").Value = "  -1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7526"
$ws.Range("E34").Value = "  +2.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.993"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("E36").Value = "  +1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.286"
$ws.Range("E37").Value = "  +2.93%  "

$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5564"
$ws.Range("E39").Value = "  +5.28%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.484"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.072"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.566"
$ws.Range("E42").Value = "  +1.77%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.736"
$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.16"
$ws.Range("E44").Value = "  -1.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1500"
$ws.Range("E45").Value = "  +1.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4763"
$ws.Range("E46").Value = "  +2.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9994"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.07"
$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.14"
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.30"
$ws.Range("E51").Value = "  +0.29%  "
